# dynamic calculation of credit capacity and financing visualization added.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (agent_id 0)
$ws.Range("D2").Value = 1.2
$ws.Range("I2").Value = 0.1
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 2
$ws.Range("N2").Value = 4

# Row 3 (agent_id 1)
$ws.Range("D3").Value = 1.2
$ws.Range("I3").Value = 0.1
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 5

# Row 4 (agent_id 2)
$ws.Range("D4").Value = 1.2
$ws.Range("I4").Value = 0.1
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 8
$ws.Range("N4").Value = 2

# Row 5 (agent_id 3)
$ws.Range("D5").Value = 1.1
$ws.Range("I5").Value = 0.1
$ws.Range("K5").Value = 6
$ws.Range("L5").Value = 9

# Row 6 (agent_id 4)
$ws.Range("D6").Value = 1.1
$ws.Range("I6").Value = 0.1
$ws.Range("K6").Value = 6
$ws.Range("L6").Value = 12

# Row 7 (agent_id 5)
$ws.Range("D7").Value = 1.1
$ws.Range("I7").Value = 0.1
$ws.Range("K7").Value = 6
$ws.Range("L7").Value = 16

# Row 8 (agent_id 6)
$ws.Range("D8").Value = 1
$ws.Range("I8").Value = 0.1
$ws.Range("K8").Value = 10
$ws.Range("L8").Value = 20

# Row 9 (agent_id 7)
$ws.Range("D9").Value = 1
$ws.Range("I9").Value = 0.1
$ws.Range("K9").Value = 10
$ws.Range("L9").Value = 25

# Row 10 (agent_id 8)
$ws.Range("D10").Value = 1
$ws.Range("I10").Value = 0.1
$ws.Range("K10").Value = 10
$ws.Range("L10").Value = 30
